$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'Datos actualizados a 22 de Marzo de 2020 a las 15:16'
$ws.Cells.Item(6, 5).Value = 24694
$ws.Cells.Item(6, 7).Value = 372
$ws.Cells.Item(6, 8).Value = 1753
$ws.Cells.Item(7, 5).Value = 26382
$ws.Cells.Item(7, 7).Value = 47
$ws.Cells.Item(7, 8).Value = 349
$ws.Cells.Item(8, 2).Value = 23937
$ws.Cells.Item(8, 3).Value = 1573
$ws.Cells.Item(8, 5).Value = 23578
$ws.Cells.Item(17, 2).Value = 2249
$ws.Cells.Item(17, 3).Value = 85
$ws.Cells.Item(17, 5).Value = 2236
$ws.Cells.Item(33, 1).Value = 'Chile'
$ws.Cells.Item(33, 2).Value = 632
$ws.Cells.Item(33, 3).Value = 95
$ws.Cells.Item(33, 4).Value = 8
$ws.Cells.Item(33, 5).Value = 623
$ws.Cells.Item(33, 6).Value = 7
$ws.Cells.Item(34, 1).Value = 'Finlandia'
$ws.Cells.Item(34, 2).Value = 626
$ws.Cells.Item(34, 3).Value = 103
$ws.Cells.Item(34, 4).Value = 10
$ws.Cells.Item(34, 5).Value = 615
$ws.Cells.Item(34, 6).Value = 12
$ws.Cells.Item(35, 1).Value = 'Tailandia'
$ws.Cells.Item(35, 2).Value = 599
$ws.Cells.Item(35, 3).Value = 188
$ws.Cells.Item(35, 4).Value = 44
$ws.Cells.Item(35, 5).Value = 554
$ws.Cells.Item(35, 6).Value = 7
$ws.Cells.Item(36, 1).Value = 'Islandia'
$ws.Cells.Item(36, 2).Value = 568
$ws.Cells.Item(36, 3).Value = 95
$ws.Cells.Item(36, 4).Value = 5
$ws.Cells.Item(36, 5).Value = 562
$ws.Cells.Item(36, 6).Value = 1
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 1
$ws.Cells.Item(37, 1).Value = 'Polonia'
$ws.Cells.Item(37, 2).Value = 563
$ws.Cells.Item(37, 3).Value = 27
$ws.Cells.Item(37, 4).Value = 13
$ws.Cells.Item(37, 5).Value = 543
$ws.Cells.Item(37, 6).Value = 3
$ws.Cells.Item(37, 7).Value = 2
$ws.Cells.Item(75, 2).Value = 115
$ws.Cells.Item(75, 3).Value = 30
$ws.Cells.Item(76, 1).Value = 'Principado de Andorra'
$ws.Cells.Item(76, 2).Value = 113
$ws.Cells.Item(76, 3).Value = 25
$ws.Cells.Item(76, 4).Value = 1
$ws.Cells.Item(76, 5).Value = 111
$ws.Cells.Item(76, 6).Value = 2
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 1
$ws.Cells.Item(77, 1).Value = 'Republica Dominicana'
$ws.Cells.Item(77, 2).Value = 112
$ws.Cells.Item(77, 3).Value = 0
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 109
$ws.Cells.Item(77, 6).Value = 0
$ws.Cells.Item(78, 1).Value = 'Marruecos'
$ws.Cells.Item(78, 2).Value = 109
$ws.Cells.Item(78, 3).Value = 13
$ws.Cells.Item(78, 4).Value = 3
$ws.Cells.Item(78, 5).Value = 103
$ws.Cells.Item(78, 6).Value = 1
$ws.Cells.Item(78, 8).Value = 3
$ws.Cells.Item(79, 1).Value = 'Vietnam'
$ws.Cells.Item(79, 2).Value = 106
$ws.Cells.Item(79, 3).Value = 12
$ws.Cells.Item(79, 4).Value = 17
$ws.Cells.Item(79, 5).Value = 89
$ws.Cells.Item(79, 6).Value = 2
$ws.Cells.Item(79, 8).Value = 0
$ws.Cells.Item(80, 1).Value = 'Lituania'
$ws.Cells.Item(80, 2).Value = 105
$ws.Cells.Item(80, 3).Value = 6
$ws.Cells.Item(80, 5).Value = 103
$ws.Cells.Item(80, 6).Value = 1
$ws.Cells.Item(80, 8).Value = 1
$ws.Cells.Item(81, 1).Value = 'Jordania'
$ws.Cells.Item(81, 2).Value = 100
$ws.Cells.Item(81, 3).Value = 0
$ws.Cells.Item(81, 4).Value = 1
$ws.Cells.Item(81, 5).Value = 99
$ws.Cells.Item(81, 6).Value = 0
$ws.Cells.Item(81, 8).Value = 0
$ws.Cells.Item(82, 1).Value = 'Bosnia y Herzegovina'
$ws.Cells.Item(82, 2).Value = 94
$ws.Cells.Item(82, 3).Value = 1
$ws.Cells.Item(82, 4).Value = 2
$ws.Cells.Item(82, 5).Value = 91
$ws.Cells.Item(82, 6).Value = 1
$ws.Cells.Item(82, 8).Value = 1
$ws.Cells.Item(83, 1).Value = 'Islas Feroe'
$ws.Cells.Item(83, 2).Value = 92
$ws.Cells.Item(83, 3).Value = 0
$ws.Cells.Item(83, 4).Value = 3
$ws.Cells.Item(83, 5).Value = 89
$ws.Cells.Item(83, 6).Value = 0
$ws.Cells.Item(84, 1).Value = 'Malta'
$ws.Cells.Item(84, 2).Value = 90
$ws.Cells.Item(84, 3).Value = 17
$ws.Cells.Item(84, 5).Value = 88
$ws.Cells.Item(84, 6).Value = 1
$ws.Cells.Item(84, 8).Value = 0
$ws.Cells.Item(85, 1).Value = 'Albania'
$ws.Cells.Item(85, 2).Value = 89
$ws.Cells.Item(85, 3).Value = 13
$ws.Cells.Item(85, 4).Value = 2
$ws.Cells.Item(85, 5).Value = 85
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 2
$ws.Cells.Item(95, 1).Value = 'Azerbaiyan'
$ws.Cells.Item(95, 2).Value = 65
$ws.Cells.Item(95, 3).Value = 12
$ws.Cells.Item(95, 4).Value = 11
$ws.Cells.Item(95, 5).Value = 53
$ws.Cells.Item(95, 8).Value = 1
$ws.Cells.Item(96, 1).Value = 'Estado de Palestina'
$ws.Cells.Item(96, 2).Value = 59
$ws.Cells.Item(96, 3).Value = 6
$ws.Cells.Item(96, 4).Value = 17
$ws.Cells.Item(96, 5).Value = 42
$ws.Cells.Item(97, 1).Value = 'Kazajistan'
$ws.Cells.Item(97, 2).Value = 57
$ws.Cells.Item(97, 3).Value = 3
$ws.Cells.Item(97, 5).Value = 57
$ws.Cells.Item(97, 6).Value = 0
$ws.Cells.Item(97, 8).Value = 0
$ws.Cells.Item(98, 1).Value = 'Guadalupe'
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 55
$ws.Cells.Item(98, 6).Value = 4
$ws.Cells.Item(98, 8).Value = 1
$ws.Cells.Item(99, 1).Value = 'Senegal'
$ws.Cells.Item(99, 2).Value = 56
$ws.Cells.Item(99, 3).Value = 0
$ws.Cells.Item(99, 4).Value = 5
$ws.Cells.Item(99, 5).Value = 51
$ws.Cells.Item(100, 1).Value = 'Oman'
$ws.Cells.Item(100, 2).Value = 55
$ws.Cells.Item(100, 3).Value = 3
$ws.Cells.Item(100, 4).Value = 17
$ws.Cells.Item(100, 5).Value = 38
$ws.Cells.Item(100, 6).Value = 0
$ws.Cells.Item(101, 1).Value = 'Georgia'
$ws.Cells.Item(101, 2).Value = 54
$ws.Cells.Item(101, 3).Value = 5
$ws.Cells.Item(101, 4).Value = 1
$ws.Cells.Item(101, 6).Value = 1
$ws.Cells.Item(102, 1).Value = 'Camboya'
$ws.Cells.Item(102, 4).Value = 2
$ws.Cells.Item(102, 5).Value = 51
$ws.Cells.Item(102, 8).Value = 0
$ws.Cells.Item(115, 1).Value = 'Consejo Danes para los Refugiados'
$ws.Cells.Item(115, 3).Value = 0
$ws.Cells.Item(116, 1).Value = 'Puerto Rico'
$ws.Cells.Item(116, 3).Value = 2
$ws.Cells.Item(120, 1).Value = 'Bolivia'
$ws.Cells.Item(120, 4).Value = 0
$ws.Cells.Item(120, 5).Value = 20
$ws.Cells.Item(121, 1).Value = 'Macao'
$ws.Cells.Item(121, 2).Value = 20
$ws.Cells.Item(121, 3).Value = 1
$ws.Cells.Item(121, 4).Value = 10
$ws.Cells.Item(121, 5).Value = 10
$ws.Cells.Item(177, 1).Value = 'Nepal'
$ws.Cells.Item(177, 3).Value = 1
$ws.Cells.Item(177, 4).Value = 1
$ws.Cells.Item(177, 8).Value = 0
$ws.Cells.Item(178, 1).Value = 'Sudan'
$ws.Cells.Item(178, 2).Value = 2
$ws.Cells.Item(178, 8).Value = 1
$ws.Cells.Item(191, 1).Value = 'San Martin (Parte Holandesa)'
$ws.Cells.Item(191, 4).Value = 0
$ws.Cells.Item(191, 5).Value = 1
